$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Cells.Item(2, 4).Value = "26.039.41"
$ws.Cells.Item(2, 5).Value = "  +0.39%  "
$ws.Cells.Item(3, 4).Value = "1.645.33"
$ws.Cells.Item(3, 5).Value = "  +0.45%  "
$ws.Cells.Item(4, 5).Value = "  +0.61%  "
Set-TextValue 5 4 "216.05"
$ws.Cells.Item(5, 5).Value = "  +0.69%  "
$ws.Cells.Item(6, 5).Value = "  +0.45%  "
$ws.Cells.Item(7, 5).Value = "  +0.55%  "
$ws.Cells.Item(8, 5).Value = "  +0.68%  "
$ws.Cells.Item(9, 5).Value = "  +0.60%  "
Set-TextValue 10 4 "19.58"
$ws.Cells.Item(10, 5).Value = "  +0.00%  "
$ws.Cells.Item(11, 5).Value = "  +0.47%  "
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "1.740.75"
$ws.Cells.Item(12, 5).Value = "  +6.83%  "
$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue 13 4 "4.27"
$ws.Cells.Item(13, 5).Value = "  +0.70%  "
$ws.Cells.Item(14, 5).Value = "  +0.37%  "
$ws.Cells.Item(15, 5).Value = "  +1.64%  "
$ws.Cells.Item(16, 5).Value = "  +0.79%  "
$ws.Cells.Item(17, 4).Value = "26.059.05"
$ws.Cells.Item(17, 5).Value = "  +0.36%  "
Set-TextValue 19 4 "194.62"
$ws.Cells.Item(19, 5).Value = "  +0.68%  "
$ws.Cells.Item(20, 5).Value = "  -0.45%  "
Set-TextValue 21 4 "9.94"
$ws.Cells.Item(21, 5).Value = "  +0.27%  "
$ws.Cells.Item(22, 5).Value = "  -0.82%  "
Set-TextValue 23 4 "0.132"
$ws.Cells.Item(23, 5).Value = "  +4.96%  "
Set-TextValue 24 4 "1.79"
$ws.Cells.Item(24, 5).Value = "  -0.21%  "
$ws.Cells.Item(25, 5).Value = "  +0.49%  "
Set-TextValue 26 4 "143.62"
$ws.Cells.Item(26, 5).Value = "  -0.27%  "
$ws.Cells.Item(27, 5).Value = "  +0.58%  "
$ws.Cells.Item(28, 5).Value = "  +0.57%  "
Set-TextValue 29 4 "1.24"
$ws.Cells.Item(29, 5).Value = "  +0.50%  "
$ws.Cells.Item(30, 5).Value = "  -0.55%  "
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 31 4 "3.26"
$ws.Cells.Item(31, 5).Value = "  +1.58%  "
$ws.Cells.Item(32, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue 32 4 "3.29"
$ws.Cells.Item(32, 5).Value = "  -0.17%  "
Set-TextValue 33 4 "1.54"
$ws.Cells.Item(33, 5).Value = "  -0.15%  "
$ws.Cells.Item(34, 5).Value = "  +1.37%  "
Set-TextValue 35 4 "0.906"
$ws.Cells.Item(35, 5).Value = "  +0.42%  "
$ws.Cells.Item(36, 4).Value = "1.131.54"
$ws.Cells.Item(36, 5).Value = "  -0.39%  "
$ws.Cells.Item(37, 5).Value = "  -0.75%  "
$ws.Cells.Item(38, 5).Value = "  +0.43%  "
$ws.Cells.Item(39, 5).Value = "  +0.19%  "
$ws.Cells.Item(40, 5).Value = "  +0.70%  "
Set-TextValue 41 4 "98.97"
$ws.Cells.Item(41, 5).Value = "  -0.41%  "
Set-TextValue 42 4 "0.798"
$ws.Cells.Item(42, 5).Value = "  +0.12%  "
$ws.Cells.Item(43, 5).Value = "  +1.62%  "
Set-TextValue 44 4 "56.56"
$ws.Cells.Item(44, 5).Value = "  +0.18%  "
$ws.Cells.Item(45, 2).Value = "RenderToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 45 4 "1.49"
$ws.Cells.Item(45, 5).Value = "  +2.91%  "
$ws.Cells.Item(46, 2).Value = "Cronos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue 46 4 "0.0523"
$ws.Cells.Item(46, 5).Value = "  -1.16%  "
Set-TextValue 47 4 "7.80"
$ws.Cells.Item(47, 5).Value = "  +1.94%  "
$ws.Cells.Item(48, 5).Value = "  +0.04%  "
$ws.Cells.Item(49, 5).Value = "  +0.33%  "
Set-TextValue 50 4 "0.0951"
$ws.Cells.Item(50, 5).Value = "  -1.00%  "
Set-TextValue 51 4 "1.18"
